# Updates the crypto price/volume snapshot table on the active sheet:
#  - refreshed prices in column D (31 rows)
#  - the "Hora" marker in column G goes from 3 -> 4 for every data row (2-51)
#  - rows 7 & 8 swap: FTXToken now ranks above KuCoinToken, each keeping
#    its own coin name, link, price and volume label

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell='D2'; Value='282.48'}
    @{Cell='G2'; Value='4'}
    @{Cell='D3'; Value='20.43'}
    @{Cell='G3'; Value='4'}
    @{Cell='D4'; Value='6.204'}
    @{Cell='G4'; Value='4'}
    @{Cell='D5'; Value='0.06169'}
    @{Cell='G5'; Value='4'}
    @{Cell='D6'; Value='3.584'}
    @{Cell='G6'; Value='4'}
    @{Cell='B7'; Value='FTXToken'}
    @{Cell='C7'; Value='https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'}
    @{Cell='D7'; Value='1.507'}
    @{Cell='E7'; Value='6FTXTokenFTT'}
    @{Cell='G7'; Value='4'}
    @{Cell='B8'; Value='KuCoinToken'}
    @{Cell='C8'; Value='https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'}
    @{Cell='D8'; Value='6.549'}
    @{Cell='E8'; Value='7KuCoinTokenKCS'}
    @{Cell='G8'; Value='4'}
    @{Cell='D9'; Value='0.8196'}
    @{Cell='G9'; Value='4'}
    @{Cell='D10'; Value='0.01379'}
    @{Cell='G10'; Value='4'}
    @{Cell='D11'; Value='0.1640'}
    @{Cell='G11'; Value='4'}
    @{Cell='D12'; Value='0.08493'}
    @{Cell='G12'; Value='4'}
    @{Cell='D13'; Value='0.03479'}
    @{Cell='G13'; Value='4'}
    @{Cell='D14'; Value='0.03216'}
    @{Cell='G14'; Value='4'}
    @{Cell='D15'; Value='0.09158'}
    @{Cell='G15'; Value='4'}
    @{Cell='D16'; Value='3.718'}
    @{Cell='G16'; Value='4'}
    @{Cell='D17'; Value='0.001649'}
    @{Cell='G17'; Value='4'}
    @{Cell='D18'; Value='0.04710'}
    @{Cell='G18'; Value='4'}
    @{Cell='D19'; Value='0.006477'}
    @{Cell='G19'; Value='4'}
    @{Cell='D20'; Value='0.006166'}
    @{Cell='G20'; Value='4'}
    @{Cell='G21'; Value='4'}
    @{Cell='G22'; Value='4'}
    @{Cell='D23'; Value='3.833'}
    @{Cell='G23'; Value='4'}
    @{Cell='G24'; Value='4'}
    @{Cell='D25'; Value='0.3350'}
    @{Cell='G25'; Value='4'}
    @{Cell='G26'; Value='4'}
    @{Cell='G27'; Value='4'}
    @{Cell='G28'; Value='4'}
    @{Cell='G29'; Value='4'}
    @{Cell='G30'; Value='4'}
    @{Cell='G31'; Value='4'}
    @{Cell='G32'; Value='4'}
    @{Cell='G33'; Value='4'}
    @{Cell='G34'; Value='4'}
    @{Cell='G35'; Value='4'}
    @{Cell='G36'; Value='4'}
    @{Cell='G37'; Value='4'}
    @{Cell='G38'; Value='4'}
    @{Cell='G39'; Value='4'}
    @{Cell='D40'; Value='0.04675'}
    @{Cell='G40'; Value='4'}
    @{Cell='D41'; Value='0.007196'}
    @{Cell='G41'; Value='4'}
    @{Cell='G42'; Value='4'}
    @{Cell='D43'; Value='0.1097'}
    @{Cell='G43'; Value='4'}
    @{Cell='D44'; Value='0.01145'}
    @{Cell='G44'; Value='4'}
    @{Cell='D45'; Value='0.00006820'}
    @{Cell='G45'; Value='4'}
    @{Cell='G46'; Value='4'}
    @{Cell='D47'; Value='1.103'}
    @{Cell='G47'; Value='4'}
    @{Cell='D48'; Value='0.002919'}
    @{Cell='G48'; Value='4'}
    @{Cell='G49'; Value='4'}
    @{Cell='G50'; Value='4'}
    @{Cell='G51'; Value='4'}
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $val = $u.Value

    if ($val -match '^-?\d+(\.\d+)?$') {
        # The source sheet stores these numeric-looking entries (prices,
        # the hour code) as literal text, not real numbers. Force text
        # formatting before the write so Excel doesn't silently coerce
        # them into numbers (which would, e.g., drop trailing zeros).
        $cell.NumberFormat = "@"
    }

    $cell.Value = $val
}
